$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fuel dist urban")

# Updated RD distribution emission factors (column C) for algae CAP and HTL
# pathways; columns D:G are `=C#` formulas and recompute automatically.
$ws.Range("C4").Value = 5093.7824552240209
$ws.Range("C5").Value = 5023.1958652604862
$ws.Range("C6").Value = 41.657173986229949
$ws.Range("C7").Value = 527.78573845908159
$ws.Range("C8").Value = 4453.7529528151745
$ws.Range("C9").Value = 0.09869537321209304
$ws.Range("C10").Value = 1.1883133447180589
$ws.Range("C11").Value = 0.60557378989801824
$ws.Range("C12").Value = 1.3472060329660924
$ws.Range("C13").Value = 0.04077408879999693
$ws.Range("C14").Value = 0.02866035992117464
$ws.Range("C15").Value = 0.024510023670013949
$ws.Range("C16").Value = 0.0033183046022218594
$ws.Range("C17").Value = 0.015425727512154881
$ws.Range("C18").Value = 0.48940655600260441
$ws.Range("C19").Value = 0.0058402637809921217
$ws.Range("C20").Value = 387.85040331686662
$ws.Range("C21").Value = 0.26551552589021865
$ws.Range("C22").Value = 0.17288665151109295
$ws.Range("C23").Value = 0.21330355735173617
$ws.Range("C24").Value = 0.0094740252492059081
$ws.Range("C25").Value = 0.0050442745092267195
$ws.Range("C26").Value = 0.0052423279836688944
$ws.Range("C27").Value = 0.00054229333446464761
$ws.Range("C28").Value = 0.0018571109385136862

# Recalculate so the dependent =C# formulas in D:G pick up the new values.
$excel.Calculate()

# Move the active sheet/selection from "Fuel specs" to "Fuel dist urban",
# matching the selection left behind in the saved file.
$ws.Activate()
$ws.Range("C10").Select()
